# Change column header for due date
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 currently reads "INITIAL COMPLIANCE DUE DATE"; rename to "CURRENT COMPLIANCE DUE DATE"
$ws.Range("B1").Value = "CURRENT COMPLIANCE DUE DATE"

# Update the active selection to B1 to match the author's final selection state
$ws.Range("B1").Select()
